# feat: add 2022-Q3 data
#
# - Inserts a new worksheet "2022-Q3" right after "总计" and before the
#   existing "2022-Q2" sheet (pushing 2022-Q2 / 2022-Q1 / 2021-Q4 back by one).
# - Populates "2022-Q3" with the fund holdings for that quarter.
# - Updates the "总计" (summary) sheet with a new row for 2022-Q3 and shifts
#   the existing rows down.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned before "2022-Q2".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# NOTE: inserting a sheet shifts everyone after it, and previously grabbed
# worksheet references become stale (they are positional). Re-resolve the
# sheets we still need by name after the insert.
$summary = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

# Bring over the header row (B1:H1) plus its formatting from the 2022-Q2
# sheet - every quarterly sheet shares an identical header.
$q2.Range("B1:H1").Copy($q3.Range("B1:H1"))

# Bring over the bold/bordered style used for the numeric index column (A)
# by copying a single already-styled cell.
$q2.Range("A2").Copy($q3.Range("A2"))
$q2.Range("A2").Copy($q3.Range("A3"))
$q2.Range("A2").Copy($q3.Range("A4"))
$q2.Range("A2").Copy($q3.Range("A5"))

$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(5,1).Value = 3

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2: 013166 / 东兴宸祥量化混合A
Set-TextValue $q3.Cells.Item(2,2) "013166"
$q3.Cells.Item(2,3).Value = "东兴宸祥量化混合A"
Set-TextValue $q3.Cells.Item(2,4) "0.38"
Set-TextValue $q3.Cells.Item(2,5) "93.87"
Set-TextValue $q3.Cells.Item(2,6) "1.16"
Set-TextValue $q3.Cells.Item(2,7) "0.0044"
$q3.Cells.Item(2,8).Value = 9

# Row 3: 009327 / 东兴兴晟混合A
Set-TextValue $q3.Cells.Item(3,2) "009327"
$q3.Cells.Item(3,3).Value = "东兴兴晟混合A"
Set-TextValue $q3.Cells.Item(3,4) "0.38"
Set-TextValue $q3.Cells.Item(3,5) "79.70"
Set-TextValue $q3.Cells.Item(3,6) "1.09"
Set-TextValue $q3.Cells.Item(3,7) "0.0041"
$q3.Cells.Item(3,8).Value = 7

# Row 4: 013167 / 东兴宸祥量化混合C
Set-TextValue $q3.Cells.Item(4,2) "013167"
$q3.Cells.Item(4,3).Value = "东兴宸祥量化混合C"
Set-TextValue $q3.Cells.Item(4,4) "0.08"
Set-TextValue $q3.Cells.Item(4,5) "93.87"
Set-TextValue $q3.Cells.Item(4,6) "1.16"
Set-TextValue $q3.Cells.Item(4,7) "0.0009"
$q3.Cells.Item(4,8).Value = 9

# Row 5: 009328 / 东兴兴晟混合C
Set-TextValue $q3.Cells.Item(5,2) "009328"
$q3.Cells.Item(5,3).Value = "东兴兴晟混合C"
Set-TextValue $q3.Cells.Item(5,4) "0.07"
Set-TextValue $q3.Cells.Item(5,5) "79.70"
Set-TextValue $q3.Cells.Item(5,6) "1.09"
Set-TextValue $q3.Cells.Item(5,7) "0.0008"
$q3.Cells.Item(5,8).Value = 7

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a 2022-Q3 row and shift the
#    pre-existing rows (2022-Q2, 2022-Q1, 2021-Q4) down by one.
# ---------------------------------------------------------------------------

# Row 5 (2021-Q4) is brand new - bring over the bold/bordered index-column
# style from row 4 before filling in the values.
$summary.Range("A4").Copy($summary.Range("A5"))

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q4"
$summary.Cells.Item(5,3).Value = 6
$summary.Cells.Item(5,4).Value = 0.33

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q1"
$summary.Cells.Item(4,3).Value = 1
$summary.Cells.Item(4,4).Value = 0.29

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 6
$summary.Cells.Item(3,4).Value = 0.29

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 4
$summary.Cells.Item(2,4).Value = 0.01
